# Updated cryptos list (refreshed Price / Volume(1h) columns), matching the
# "Updated cryptos list ... with GitHub Actions" scheduled data refresh.
# Rows 30/31 also swap identity (PancakeSwap <-> InternetComputer(DFINITY))
# because the source ranking reordered between runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / already-non-numeric values: direct assignment is safe
$ws.Range("D2").Value = '65.188.04'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '3.400.13'
$ws.Range("E3").Value = '  -3.61%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  -2.99%  '
$ws.Range("E6").Value = '  -5.43%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.400.54'
$ws.Range("E8").Value = '  -3.57%  '
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("E10").Value = '  -8.20%  '
$ws.Range("E11").Value = '  -11.42%  '
$ws.Range("E12").Value = '  -7.99%  '
$ws.Range("D13").Value = '3.981.42'
$ws.Range("E13").Value = '  -3.58%  '
$ws.Range("E14").Value = '  -11.31%  '
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").Value = '3.415.35'
$ws.Range("E16").Value = '  -2.65%  '
$ws.Range("D17").Value = '65.187.56'
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("E18").Value = '  -10.29%  '
$ws.Range("E19").Value = '  -12.12%  '
$ws.Range("E20").Value = '  -5.94%  '
$ws.Range("E21").Value = '  -5.10%  '
$ws.Range("E22").Value = '  -8.15%  '
$ws.Range("E23").Value = '  -8.34%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -6.73%  '
$ws.Range("D26").Value = '3.539.16'
$ws.Range("E26").Value = '  -3.58%  '
$ws.Range("E27").Value = '  -10.52%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("E29").Value = '  -9.84%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E30").Value = '  -10.09%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("E31").Value = '  -10.37%  '
$ws.Range("D32").Value = '3.409.82'
$ws.Range("E32").Value = '  -3.34%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  -7.52%  '
$ws.Range("E35").Value = '  -6.98%  '
$ws.Range("E36").Value = '  -3.79%  '
$ws.Range("E37").Value = '  -11.58%  '
$ws.Range("E38").Value = '  -11.78%  '
$ws.Range("E39").Value = '  -8.01%  '
$ws.Range("E40").Value = '  -11.92%  '
$ws.Range("E41").Value = '  -8.48%  '
$ws.Range("E42").Value = '  -6.24%  '
$ws.Range("E43").Value = '  -4.04%  '
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("E45").Value = '  -14.39%  '
$ws.Range("E46").Value = '  -10.83%  '
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("E48").Value = '  -2.66%  '
$ws.Range("E49").Value = '  -8.88%  '
$ws.Range("E50").Value = '  -15.29%  '
$ws.Range("D51").Value = '2.169.41'
$ws.Range("E51").Value = '  -7.95%  '

# Numeric-looking price values must be forced to Text so Excel doesn't
# normalize/round them (e.g. '13.60' -> 13.6). Force text format, assign,
# then restore the default style so no stray per-cell format lingers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.119'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.372'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000177'
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '380.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.549'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000104'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.02'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.142'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.63'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '168.75'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.68'
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0752'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.805'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.46'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.02'
$ws.Range("D50").Style = "Normal"
